# Actualización automática 2025-12-07 17:54:30
#
# In the "VENTA MENSUAL" sheet, the duplicate row for
# "PAUTA ASTUDILLO JULIO HERNAN" (row 18) is removed. All rows below it
# shift up by one, and the totals row (previously row 22, now row 21) is
# updated to reflect the removal of that row's D and G values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Remove row 18 (duplicate "PAUTA ASTUDILLO JULIO HERNAN" entry); all
# subsequent rows shift up by one automatically.
$ws.Rows.Item(18).Delete()

# The totals row was row 22 (D=17549.57, G=1000) and is now row 21.
# Since these totals are static values (not formulas), update them to
# remove the contribution of the deleted row (D -= 326.73, G -= 1000).
$ws.Range("D21").Value = 17222.84
$ws.Range("G21").Value = 0
